$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NextYear")
$ws.Select()

# Data computed by the NN (NextYear) model - test-9, fills Z7:AA26
$values = @(
    @(7, 560.50998950074211, 1275.600167235162),
    @(8, 522.05038851027052, 1245.00450099573),
    @(9, 510.54175828544771, 1224.0457560982541),
    @(10, 594.83201247958345, 1250.881080885373),
    @(11, 533.71593172848497, 1233.518200058767),
    @(12, 490.57063055947663, 1424.9429228418469),
    @(13, 500.95261468811299, 1317.0755869606271),
    @(14, 526.10675567673206, 1268.094857237646),
    @(15, 493.94392799948582, 1352.9965724416379),
    @(16, 487.81024884027039, 1309.1864162922529),
    @(17, 501.7510623320257, 1256.8473536241429),
    @(18, 601.22116389570601, 1327.9264803926901),
    @(19, 565.4341859150062, 1135.0631737094229),
    @(20, 608.49863593354053, 1204.007012169282),
    @(21, 521.39625719769197, 1335.897831061056),
    @(22, 524.26255120630321, 1399.736875125539),
    @(23, 544.10427197046988, 1292.224007560751),
    @(24, 518.69560857377508, 1217.131283365979),
    @(25, 496.69156082074068, 1235.834797935581),
    @(26, 471.4935830244371, 1343.899372083657)
)

foreach ($row in $values) {
    $r = $row[0]
    $ws.Range("Z$r").Value = $row[1]
    $ws.Range("AA$r").Value = $row[2]
}

# The cells Z7:AA26 (and the blank/summary rows right below, 27-29) previously
# used a placeholder scientific number format (0.00E+00) reserved for
# not-yet-computed results. Now that real data is present, restyle the whole
# block to match the other data columns on the sheet (0.00, centered).
$dataRange = $ws.Range("Z7:AA29")
$dataRange.NumberFormat = "0.00"
$dataRange.HorizontalAlignment = -4108

# Update the active selection to V7 (matches the author's final cursor position)
[void]$ws.Range("V7").Select()
